$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list on Thu Mar  7 10:36:44 UTC 2024 with GitHub Actions

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "66.688.38"
$ws.Cells.Item(2, 5).Value = "  +0.36%  "

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "3.787.45"
$ws.Cells.Item(3, 5).Value = "  -1.01%  "

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.00"
$ws.Cells.Item(4, 5).Value = "  +0.16%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "434.29"
$ws.Cells.Item(5, 5).Value = "  +2.35%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "139.74"
$ws.Cells.Item(6, 5).Value = "  +6.99%  "

$ws.Cells.Item(7, 5).Value = "  +1.50%  "

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "1.00"
$ws.Cells.Item(8, 5).Value = "  +0.04%  "

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.735"
$ws.Cells.Item(9, 5).Value = "  +0.67%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.152"
$ws.Cells.Item(10, 5).Value = "  -9.46%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0000314"
$ws.Cells.Item(11, 5).Value = "  -14.28%  "

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "42.96"
$ws.Cells.Item(12, 5).Value = "  +5.01%  "

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "10.40"
$ws.Cells.Item(13, 5).Value = "  +3.22%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "4.398.30"
$ws.Cells.Item(14, 5).Value = "  -0.71%  "

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "14.86"
$ws.Cells.Item(15, 5).Value = "  -4.13%  "

$ws.Cells.Item(16, 2).Value = "TRON"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "0.137"
$ws.Cells.Item(16, 5).Value = "  -0.49%  "

$ws.Cells.Item(17, 2).Value = "WrappedEther"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "3.761.18"
$ws.Cells.Item(17, 5).Value = "  -1.74%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "19.92"
$ws.Cells.Item(18, 5).Value = "  +1.58%  "

$ws.Cells.Item(19, 5).Value = "  +7.18%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "66.859.00"
$ws.Cells.Item(20, 5).Value = "  -0.02%  "

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "409.80"
$ws.Cells.Item(21, 5).Value = "  -0.53%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "14.76"
$ws.Cells.Item(22, 5).Value = "  +2.11%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "3.24"
$ws.Cells.Item(23, 5).Value = "  +6.76%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "85.59"
$ws.Cells.Item(24, 5).Value = "  +0.07%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "36.85"
$ws.Cells.Item(25, 5).Value = "  -0.31%  "

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "3.38"
$ws.Cells.Item(26, 5).Value = "  +4.84%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "9.87"
$ws.Cells.Item(27, 5).Value = "  +37.84%  "

$ws.Cells.Item(28, 5).Value = "  -2.36%  "

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "9.77"
$ws.Cells.Item(29, 5).Value = "  +2.89%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "725.70"
$ws.Cells.Item(30, 5).Value = "  +6.06%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "13.75"
$ws.Cells.Item(31, 5).Value = "  +10.28%  "

$ws.Cells.Item(32, 5).Value = "  +9.68%  "

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "2.74"
$ws.Cells.Item(33, 5).Value = "  +0.00%  "

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "42.66"

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "1.00"
$ws.Cells.Item(35, 5).Value = "  -0.04%  "

$ws.Cells.Item(36, 5).Value = "  +0.85%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "5.63"
$ws.Cells.Item(37, 5).Value = "  +24.97%  "

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "56.16"
$ws.Cells.Item(38, 5).Value = "  +2.34%  "

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.0478"
$ws.Cells.Item(39, 5).Value = "  +4.61%  "

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "2.76"
$ws.Cells.Item(40, 5).Value = "  +42.65%  "

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "2.98"
$ws.Cells.Item(41, 5).Value = "  -5.60%  "

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.141"
$ws.Cells.Item(42, 5).Value = "  +3.44%  "

$ws.Cells.Item(43, 2).Value = "PEPE"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.0₃0670"
$ws.Cells.Item(43, 5).Value = "  -15.66%  "

$ws.Cells.Item(44, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "1.00"
$ws.Cells.Item(44, 5).Value = "  -0.33%  "

$ws.Cells.Item(45, 2).Value = "ApeXProtocol"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "3.33"
$ws.Cells.Item(45, 5).Value = "  +5.79%  "

$ws.Cells.Item(46, 2).Value = "TheGraph"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.332"
$ws.Cells.Item(46, 5).Value = "  +13.19%  "

$ws.Cells.Item(47, 5).Value = "  +5.66%  "

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "3.31"
$ws.Cells.Item(48, 5).Value = "  +0.54%  "

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "2.08"
$ws.Cells.Item(49, 5).Value = "  +0.32%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "142.05"
$ws.Cells.Item(50, 5).Value = "  -4.65%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "2.81"
$ws.Cells.Item(51, 5).Value = "  +1.07%  "

